# Fill in extra "fake data" answer rows (2-5) on the "antwoorden" sheet,
# restyle the answer cells (smaller wrapped font) and update the
# workbook's saved absolute path / selection to match the re-export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("antwoorden")

# ---------------------------------------------------------------------
# 1. Row 2 (B2:Z2) gets a new set of fake answers (A2 already holds "A"
#    and keeps it).
# ---------------------------------------------------------------------
$row2 = @("A","A","D","D","A","C","C","D","B","B","D","D","A","C","A","D","D","D","B","B","B","D","C","B","D")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

# ---------------------------------------------------------------------
# 2. Three new fake-data rows (3, 4 and 5), columns A:Z.
# ---------------------------------------------------------------------
$row3 = @("B","C","B","B","C","A","D","C","A","B","C","C","B","D","A","C","D","D","C","B","B","A","D","C","B","B")
$row4 = @("C","B","D","C","C","C","C","C","C","D","A","A","D","C","D","D","A","A","A","A","B","A","C","C","A","A")
$row5 = @("D","B","D","A","A","D","C","D","D","A","D","C","B","B","B","C","C","D","A","A","C","D","A","A","D","A")

$newRows = @($row3, $row4, $row5)
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $values = $newRows[$r]
    $rowIndex = $r + 3
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $values[$c]
    }
}

# ---------------------------------------------------------------------
# 3. Restyle the answer block: smaller "Liberation Sans" font, wrapped
#    and vertically centred text. Update A2 first (the sole user of the
#    existing "left aligned" style) so the style definition itself is
#    updated in place, then fan that exact style out to the rest of the
#    answer block (B2:Z5) with a format-only paste.
# ---------------------------------------------------------------------
$anchor = $ws.Range("A2")
$anchor.Font.Size = 10
$anchor.Font.Name = "Liberation Sans"
$anchor.VerticalAlignment = -4108   # xlCenter
$anchor.WrapText = $true

$anchor.Copy() | Out-Null
$ws.Range("A2:Z5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Dimension grows to A1:Z5 automatically with the new data; reflect
#    the matching selection in the saved sheet view.
# ---------------------------------------------------------------------
$ws.Range("A1:Z5").Select() | Out-Null

